# ENH: re-run all analyses
# Update epl (column J) values in the supplementary results table with
# refreshed figures from the re-run analysis notebook. Cells are plain
# text (stored as inline strings originally) so we prefix each literal
# with an apostrophe to force Excel to keep storing them as text rather
# than coercing the numeric-looking strings into Number cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J3").Value = "'0.014"
$ws.Range("J4").Value = "'0.017"
$ws.Range("J5").Value = "'0.119"
$ws.Range("J6").Value = "'0.126"
$ws.Range("J7").Value = "'0.139"
$ws.Range("J9").Value = "'0.035"
$ws.Range("J10").Value = "'0.033"
$ws.Range("J11").Value = "'0.104"
$ws.Range("J12").Value = "'0.165"
$ws.Range("J13").Value = "'0.200"
$ws.Range("J15").Value = "'0.038"
$ws.Range("J16").Value = "'0.046"
$ws.Range("J17").Value = "'0.166"
$ws.Range("J18").Value = "'0.150"
$ws.Range("J19").Value = "'0.225"
$ws.Range("J21").Value = "'0.048"
$ws.Range("J22").Value = "'0.064"
$ws.Range("J23").Value = "'0.166"
$ws.Range("J24").Value = "'0.164"
$ws.Range("J25").Value = "'0.241"
$ws.Range("J27").Value = "'0.094"
$ws.Range("J28").Value = "'0.116"
$ws.Range("J29").Value = "'0.300"
$ws.Range("J30").Value = "'0.253"
$ws.Range("J31").Value = "'0.280"
$ws.Range("J33").Value = "'0.169"
$ws.Range("J34").Value = "'0.219"
$ws.Range("J35").Value = "'0.453"
$ws.Range("J36").Value = "'0.445"
$ws.Range("J39").Value = "'0.225"
$ws.Range("J40").Value = "'0.299"
$ws.Range("J41").Value = "'0.541"
$ws.Range("J42").Value = "'0.626"
$ws.Range("J43").Value = "'0.652"
$ws.Range("J45").Value = "'0.318"
$ws.Range("J46").Value = "'0.381"
$ws.Range("J47").Value = "'0.742"
$ws.Range("J48").Value = "'0.722"
$ws.Range("J49").Value = "'0.743"
